# Horarios Línea 141 update - 2026-02-01 run at 05:01:49 (scrape ts 05:01:48)
$wb = $excel.ActiveWorkbook

$updated = "Última actualización: 05:01:49"

# ---------------------------------------------------------------
# Sheet "LP1912": rows 6-17 get refreshed data, a new row 18 is added
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = $updated
$ws1.Range("A3").Value = "Total filas: 13"

$sheet1Data = @(
    @("05:01:48", "05:16", "17_ROMERO",     15, "LP1912"),
    @("05:01:48", "05:22", "23_HERNANDEZ",  21, "LP1912"),
    @("05:01:48", "05:44", "14_ABASTO",     43, "LP1912"),
    @("05:01:48", "05:46", "17_ROMERO",     45, "LP1912"),
    @("05:01:48", "06:01", "16_SANTA ANA",  60, "LP1912"),
    @("05:01:48", "06:09", "10_OLMOS",      68, "LP1912"),
    @("05:01:48", "06:15", "215A_EL PATO",  74, "LP1912"),
    @("05:01:48", "06:30", "23_HERNANDEZ",  89, "LP1912"),
    @("05:01:48", "06:34", "11_ETCHEVERRY", 93, "LP1912"),
    @("05:01:48", "06:38", "17X38_ROMERO",  97, "LP1912"),
    @("05:01:48", "06:40", "16_SANTA ANA",  99, "LP1912"),
    @("05:01:48", "06:56", "215A_EL PATO", 115, "LP1912"),
    @("05:01:48", "06:59", "225_GOMEZ",    118, "LP1912")
)

$row = 6
foreach ($rec in $sheet1Data) {
    $ws1.Cells.Item($row, 1).Value = $rec[0]
    $ws1.Cells.Item($row, 2).Value = $rec[1]
    $ws1.Cells.Item($row, 3).Value = $rec[2]
    $ws1.Cells.Item($row, 4).Value = $rec[3]
    $ws1.Cells.Item($row, 5).Value = $rec[4]
    $row = $row + 1
}

# ---------------------------------------------------------------
# Sheet "LP1912-215": row 6 refreshed, new row 7 added
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = $updated
$ws2.Range("A3").Value = "Total filas: 2"

$sheet2Data = @(
    @("05:01:48", "06:15", "215A_EL PATO",  74, "LP1912"),
    @("05:01:48", "06:56", "215A_EL PATO", 115, "LP1912")
)

$row = 6
foreach ($rec in $sheet2Data) {
    $ws2.Cells.Item($row, 1).Value = $rec[0]
    $ws2.Cells.Item($row, 2).Value = $rec[1]
    $ws2.Cells.Item($row, 3).Value = $rec[2]
    $ws2.Cells.Item($row, 4).Value = $rec[3]
    $ws2.Cells.Item($row, 5).Value = $rec[4]
    $row = $row + 1
}

# ---------------------------------------------------------------
# Sheet "6203-6173": only the timestamp refreshes, no rows change
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = $updated
